$d = $word.ActiveDocument

$pairs = @(
    @("86×50=4300", "87×54=4698"),
    @("29×95=2755", "68×71=4828"),
    @("87×83=7221", "79×35=2765"),
    @("17×96=1632", "12×27=324"),
    @("74×94=6956", "82×18=1476"),
    @("68×47=3196", "55×50=2750"),
    @("19×80=1520", "86×62=5332"),
    @("24×72=1728", "35×23=805"),
    @("29×36=1044", "25×74=1850"),
    @("59×82=4838", "36×52=1872"),
    @("71×80=5680", "30×37=1110"),
    @("92×37=3404", "21×40=840"),
    @("63×31=1953", "30×23=690"),
    @("31×38=1178", "90×93=8370"),
    @("49×40=1960", "98×11=1078"),
    @("21×50=1050", "53×82=4346"),
    @("85×26=2210", "18×38=684"),
    @("73×46=3358", "83×96=7968"),
    @("90×72=6480", "90×96=8640"),
    @("63×46=2898", "18×29=522"),
    @("51×94=4794", "79×37=2923"),
    @("63×22=1386", "15×60=900"),
    @("58×32=1856", "20×46=920"),
    @("95×43=4085", "93×51=4743"),
    @("67×88=5896", "92×20=1840")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
